# Remove the trailing "Ver no Jupiter..." / copyright footer block.
# That block consists of:
#   - a blank paragraph right after "LOQ4055: Quimica Inorgânica (Requisito fraco)"
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ..."
# The final blank paragraph (just before the page-break paragraph) is kept.
$d = $word.ActiveDocument
$paras = $d.Paragraphs

$jupiterMarker = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightMarker = "Contact: luizeleno@usp.br"

$jupiterIndex = 0
$copyrightIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $text = $paras.Item($i).Range.Text
    if ($text -like "*$jupiterMarker*") {
        $jupiterIndex = $i
    }
    if ($text -like "*$copyrightMarker*") {
        $copyrightIndex = $i
    }
}

$blankIndex = $jupiterIndex - 1

$delStart = $paras.Item($blankIndex).Range.Start
$delEnd = $paras.Item($copyrightIndex).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()
